# Digital Transformation Summary - French text revisions
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Fabrikam Inc. a fait l’objet d’une initiative complète de transformation numérique visant à améliorer l’efficacité opérationnelle, à améliorer l’expérience client et à stimuler l’innovation." "Fabrikam Inc. a mené une initiative de transformation numérique complète visant à améliorer son efficacité opérationnelle, à améliorer l’expérience client et à stimuler l’innovation."

Replace-Text "Le résumé suivant décrit les principales mises à jour et les jalons réalisés à ce jour." "Le résumé suivant décrit les principales mises à jour et les étapes majeures réalisées à ce jour."

Replace-Text "Mises à jour clés" "Mises à jour principales"

Replace-Text "Implémentation de l’infrastructure cloud" "Implémentation d’une infrastructure cloud"

Replace-Text "Migration de 80 % des applications locales vers le cloud." "Migration de 80 % des applications locales vers le cloud"

Replace-Text "Amélioration de la scalabilité et réduction des coûts informatiques de 25 %." "Amélioration de la scalabilité et réduction des coûts informatiques de 25 %"

Replace-Text "Sécurité et conformité des données améliorées avec les normes du secteur." "Amélioration de la sécurité et de la conformité des données avec les normes du secteur"

Replace-Text "Intégration d’analyses basées sur l’IA pour simplifier les processus décisionnels." "Intégration d’analyses basées sur l’IA pour simplifier les processus décisionnels"

Replace-Text "Déploiement de modèles Machine Learning pour prédire le comportement des clients et personnaliser les efforts marketing." "Déploiement de modèles Machine Learning pour prédire le comportement des clients et personnaliser les initiatives marketing"

Replace-Text "Réduction des processus manuels, ce qui entraîne une augmentation de 30 % de la productivité." "Réduction des processus manuels, entraînant une augmentation de 30 % de la productivité"

Replace-Text "Lancement d’un nouveau portail client avec des fonctionnalités en libre-service." "Lancement d’un nouveau portail client avec des fonctionnalités en libre-service"

Replace-Text "Introduction des chatbots pour le support client 24/7, ce qui réduit les temps de réponse de 50 %." "Introduction de chatbots pour un support client 24/7, réduisant les temps de réponse de 50 %"

Replace-Text "Amélioration de la satisfaction des clients de 20 % au cours de l’année dernière." "Amélioration de la satisfaction des clients de 20 % au cours de l’année dernière"

Replace-Text "Implémentation de l’automatisation des processus robotisés (RPA) pour les tâches courantes." "Implémentation de l’automatisation robotisée des processus (RPA) pour les tâches courantes"

Replace-Text "A atteint une réduction de 40 % du temps de traitement pour les opérations commerciales clés." "Réduction de 40 % du temps de traitement pour les opérations clés de l’entreprise"

Replace-Text "Réaffectation des ressources humaines à des rôles plus stratégiques au sein de l’organisation." "Réaffectation des ressources humaines à des rôles plus stratégiques au sein de l’organisation"

Replace-Text "Des programmes d’alphabétisation numérique ont été menés pour tous les employés." "Programmes de formation numérique pour tous les employés"

Replace-Text "A lancé une nouvelle plateforme d’apprentissage électronique avec des cours sur les technologies émergentes." "Nouvelle plateforme de formation en ligne avec des cours sur les technologies émergentes"

Replace-Text "Augmentation de l’engagement des employés et de l’adoption de nouveaux outils de 35 %." "Augmentation de 35 %. de l’engagement des employés et de l’adoption de nouveaux outils"

Replace-Text "Q1 2024 : Migration terminée vers l’infrastructure cloud." "T1 2024 : Migration vers l’infrastructure cloud"

Replace-Text "Q2 2024 : Plateforme d’analytique basée sur l’IA lancée." "T2 2024 : Lancement de la plateforme d’analytique basée sur l’IA"

Replace-Text "Q3 2024 : Introduction du nouveau portail client numérique." "T3 2024 : Introduction du nouveau portail client numérique"

Replace-Text "Q4 2024 : 50 % d’automatisation des processus de routine." "T4 2024 : 50 % d’automatisation des processus courants"

Replace-Text "Poursuivez l’expansion des applications IA et Machine Learning dans tous les services." "Poursuivre l’expansion des applications IA et Machine Learning dans tous les services"

Replace-Text "Améliorez davantage l’expérience client numérique avec de nouvelles fonctionnalités et services." "Améliorer davantage l’expérience client numérique avec de nouvelles fonctionnalités et de nouveaux services"

Replace-Text "Concentrez-vous sur les mesures de cybersécurité pour vous protéger contre les menaces en constante évolution." "Se concentrer sur les mesures de cybersécurité pour se protéger contre les menaces en constante évolution"

Replace-Text "Développez une stratégie numérique complète pour les cinq prochaines années." "Développer une stratégie numérique complète pour les cinq prochaines années"

Replace-Text "L’organisation reste engagée à tirer parti de la technologie pour stimuler la croissance et l’innovation futures." "L’organisation reste déterminée à tirer parti de la technologie pour stimuler la croissance et l’innovation."

# Bold the two sub-headings that should now be bold (Heading3 styled text
# that previously had explicit w:b val=0 overriding the style).
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "Automatisation des processus" -or $t -eq "Plans pour l’avenir") {
        $p.Range.Font.Bold = 1
    }
}
